# Update fortune document: rewrite the 32 fortune messages (A2:A33) with
# the revised wording, then drop the now-unused final row (old A34),
# shrinking the used range from A1:A34 to A1:A33.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "가위바위보를 한다면\r\n주먹을 내보는 건 어떨까요?"
$ws.Range("A3").Value = "연락을 망설이고 있다면\r\n오늘이 기회일 것 같아요!"
$ws.Range("A4").Value = "진심어린 사과를\r\n먼저 해보면 어떨까요?"
$ws.Range("A5").Value = "새로 알게된 인연에게\r\n관심이 생긴다면\r\n적극적으로 다가가 보세요!"
$ws.Range("A6").Value = "조금만 더 노력한다면\r\n좋은 일이 생길 거예요!"
$ws.Range("A7").Value = "미래를 위해 이번 달은\r\n절약을 해야 할 것 같아요"
$ws.Range("A8").Value = "고민 하고 있다면\r\n한번 해보는 게 어때요?"
$ws.Range("A9").Value = "감정을 숨기지 마세요!\r\n당신을 기다리고 있어요!"
$ws.Range("A10").Value = "하고 싶은 일이 있군요?\r\n망설이지 마세요\r\n잘 될 거예요!"
$ws.Range("A11").Value = "칭찬을 해보는 건 어떨까요?\r\n고래가 춤을 추길 바라면서요!"
$ws.Range("A12").Value = "솔직한 모습을 보여줄 때\r\n당신의 매력은 더욱 커져요!"
$ws.Range("A13").Value = "지금 지쳐있나요?\r\n재정비를 해보세요!"
$ws.Range("A14").Value = "오늘은 나의 언행을\r\n조심할 필요가 있어요!"
$ws.Range("A15").Value = "영감이 필요하다면\r\n영화를 보는 건 어떨까요?"
$ws.Range("A16").Value = "섣부른 결정은 피하세요!\r\n타이밍이 중요해요"
$ws.Range("A17").Value = "노력은 배신하지 않아요\r\n성실한 당신을 믿어주세요!"
$ws.Range("A18").Value = "당신은 큰 존재입니다!\r\n주변은 당신에게 의지해요"
$ws.Range("A19").Value = "웃는 얼굴이 제일 예뻐요\r\n당신의 미소가 힘이 돼요"
$ws.Range("A20").Value = "너무 걱정 말아요\r\n별 것 아닐 거예요!"
$ws.Range("A21").Value = "내일의 나는\r\n오늘의 나보다 나을 거예요"
$ws.Range("A22").Value = "당신이 믿는 신념은\r\n틀리지 않을 거예요"
$ws.Range("A23").Value = "얼마 전을 떠올려 보세요\r\n분명 계속 성장했어요!"
$ws.Range("A24").Value = "오늘 하루는 무엇이든\r\n긍정적으로 생각해보세요!"
$ws.Range("A25").Value = "무엇보다도 나 자신을\r\n사랑하는 것이 중요해요"
$ws.Range("A26").Value = "행복은\r\n거창한 것이 아니라\r\n작은 것에서 시작해요"
$ws.Range("A27").Value = "이번 방학엔\r\n나를 위한 여행을\r\n떠나보세요!"
$ws.Range("A28").Value = "오지 않은 미래를\r\n걱정하는 것보다\r\n마주한 현재에\r\n최선을 다해보세요!"
$ws.Range("A29").Value = "모든 것은\r\n나의 태도에 달렸어요!"
$ws.Range("A30").Value = "우리는 아직\r\n우리가 가치를 몰라요\r\n당신은 특별해요!"
$ws.Range("A31").Value = "익숙함에 속아\r\n소중함을 잃지 말자구요!"
$ws.Range("A32").Value = "오늘은 부모님께\r\n전화 한 통 드려보세요!"
$ws.Range("A33").Value = "한 번 뿐인 인생\r\n무엇을 망설이나요?"

$ws.Rows(34).Delete()

$ws.Range("A28").Select()

